# Updated data acquisition log with all downloaded budget documents and G.O.
#
# Original sheet had rows 1 (header), 2 (AP Finance / Volume-I-1), and 4
# (AP Finance / Volume-III-11 -- Agriculture). This edit:
#   - moves the old row 4 content up to row 3 (same data, reformatted to
#     match row 2's plain cell styling, keeping its tall 100.8 row height)
#   - adds three brand-new rows (4, 5, 6) for the newly downloaded RD
#     budget documents and the Panchayat Raj G.O. / pension scheme entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hyperlink that lived on B4 -- its target cell is being
# completely rebuilt below (B4 becomes a different URL entirely).
$ws.Range("B4").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# Row 3: Agriculture & Co-Operation ... (was row 4), ht 100.8
# ---------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)

$ws.Range("A3").Value = "Andhra Pradesh Finance Department"
$ws.Range("B3").Value = "https://apfinance.gov.in/...Bud@et24-25/documents/Volume-III-11.pdf"
$ws.Range("C3").Value = "Agriculture & Co-Operation & Food, Civil Supplies & Consumers Affairs Department"
$ws.Range("D3").Value = "2024-25"
$ws.Range("E3").Value2 = 45845
$ws.Range("F3").Value = "Download from website"
$ws.Range("G3").Value = "Scheme-wise data for Annadata Sukhibhava (current farmer scheme) and Vaddi Leni Runalu (VLR) found in 'LIST OF SCHEMES' table starting on Page 67. All figures are in Rupees Lakhs. YSR Rythu Bharosa was the previous scheme, Annadata Sukhibhava is the current one"
$ws.Range("H3").Value = """C:\Project_AP_Welfare_Dashboard\Data\Raw\Agriculture_data.pdf"""

# G3 needs the bordered+wrap style (like C3/H3), not the plain one copied from G2
$ws.Range("C2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "Scheme-wise data for Annadata Sukhibhava (current farmer scheme) and Vaddi Leni Runalu (VLR) found in 'LIST OF SCHEMES' table starting on Page 67. All figures are in Rupees Lakhs. YSR Rythu Bharosa was the previous scheme, Annadata Sukhibhava is the current one"

$ws.Hyperlinks.Add($ws.Range("B3"), "https://apfinance.gov.in/...Bud@et24-25/documents/Volume-III-11.pdf") | Out-Null
# Hyperlinks.Add() force-applies its own style; restore the table's normal hyperlink look
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Rows.Item(3).RowHeight = 100.8

# ---------------------------------------------------------------------
# Row 4: 2024-25 RD Data, ht 57.6
# ---------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Andhra Pradesh Finance Department"
$ws.Range("B4").Value = "https://apfinance.gov.in/...Bud@et25-26/documents/Volume-III-14.pdf"
$ws.Range("C4").Value = "2024-25 RD Data"
$ws.Range("D4").Value = "2024-25"
$ws.Range("E4").Value2 = 45845
$ws.Range("F4").Value = "Download from website"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "C:\Project_AP_Welfare_Dashboard\Data\Raw\2024-25 RD Data.pdf"

$ws.Hyperlinks.Add($ws.Range("B4"), "https://apfinance.gov.in/...Bud@et25-26/documents/Volume-III-14.pdf") | Out-Null
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Rows.Item(4).RowHeight = 57.6

# ---------------------------------------------------------------------
# Row 5: 2025-26 RD Data, ht 57.6
# ---------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

$ws.Range("A5").Value = "Andhra Pradesh Finance Department"
$ws.Range("B5").Value = "https://apfinance.gov.in/...Bud@et24-25/documents/Volume-III-14.pdf"
$ws.Range("C5").Value = "2025-26 RD Data"
$ws.Range("D5").Value = "2025-26"
$ws.Range("E5").Value2 = 45845
$ws.Range("F5").Value = "Download from website"
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = "C:\Project_AP_Welfare_Dashboard\Data\Raw\2025-26 RD Data.pdf"

$ws.Hyperlinks.Add($ws.Range("B5"), "https://apfinance.gov.in/...Bud@et24-25/documents/Volume-III-14.pdf") | Out-Null
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Rows.Item(5).RowHeight = 57.6

# ---------------------------------------------------------------------
# Row 6: Panchayat Raj G.O. / New Pension scheme, ht 57.6
# ---------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)

$ws.Range("A6").Value = "PANCHAYAT RAJ AND RURAL DEVELOPMENT (RD.I) DEPARTMENT"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "New_Pension_scheme"
$ws.Range("D6").Value = "2025-26"
$ws.Range("E6").Value2 = 45845
$ws.Range("F6").Value = "Download from website"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "C:\Project_AP_Welfare_Dashboard\Data\Raw\New_Pensions_Scheme.pdf"

# A6 gets the bordered+wrap style (long department name), B6 stays a plain
# bordered cell (no hyperlink font, since there is no link here)
$ws.Range("C2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "PANCHAYAT RAJ AND RURAL DEVELOPMENT (RD.I) DEPARTMENT"

$ws.Range("A2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = ""

$ws.Rows.Item(6).RowHeight = 57.6

# ---------------------------------------------------------------------
# Selection / housekeeping to match the final view state
# ---------------------------------------------------------------------
$ws.Range("A1:H6").Select()
$ws.Cells.Item(6, 7).Activate()

Write-Host "edit complete"
